# calorimetry : scripts : data load : all setup info accumulated into setup file/sheet
#
# The "targets" sheet only held a single "constants " label; that label now
# moves onto the "setup" sheet (alongside the "HL" reaction it refers to),
# and the "targets" sheet itself is removed. "enthalpies" shifts up to take
# its place and becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Move the "constants " / "HL" info that lived on "targets" onto "setup".
$setup = $wb.Worksheets.Item("setup")
$setup.Range("A4").Value = "constants "
$setup.Range("B4").Value = "HL"
[void]$setup.Range("B5").Select()

# Drop the now-redundant "targets" sheet.
$targets = $wb.Worksheets.Item("targets")
[void]$targets.Delete()

# "enthalpies" becomes the active sheet/tab.
$enthalpies = $wb.Worksheets.Item("enthalpies")
[void]$enthalpies.Activate()
[void]$enthalpies.Range("C8").Select()
